# Append two new rows (17 and 18) of order data to the sheet.
# Every column in this sheet (including numeric-looking Quantity / Cost
# values) is stored as text, so we force text entry with a leading
# apostrophe and then clear the resulting "quote prefix" formatting so the
# new cells end up with the same default style as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "183096-7"
$ws.Range("B17").Value = "Clio - Greek Yogurt Bar Strawberry"
$ws.Range("C17").Value = "'1"
$ws.Range("D17").Value = "'15.45"
$ws.Range("E17").Value = "'15.45"
$ws.Range("C17:E17").ClearFormats()

$ws.Range("A18").Value = "183090-0"
$ws.Range("B18").Value = "Clio - Greek Yogurt Bar Vanilla"
$ws.Range("C18").Value = "'1"
$ws.Range("D18").Value = "'15.45"
$ws.Range("E18").Value = "'15.45"
$ws.Range("C18:E18").ClearFormats()
